# Insert a new weekly observation row at row 74 (pushing the existing
# rows 74..218 down to 75..219), then populate the newly-inserted row
# with the new data point. The rest of the data keeps its values, just
# shifted down by one row - matching the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 74 downwards (row 74 becomes a fresh blank row).
$ws.Rows.Item(74).EntireRow.Insert()

# Populate the new row 74 with the new weekly record.
$ws.Cells.Item(74, 1).Value = 5
$ws.Cells.Item(74, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(74, 3).Value = "Maule"
$ws.Cells.Item(74, 4).Value = 44533
$ws.Cells.Item(74, 5).Value = 7
$ws.Cells.Item(74, 6).Value = 100114014
$ws.Cells.Item(74, 7).Value = "Betarraga"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 5000
$ws.Cells.Item(74, 11).Value = 500
$ws.Cells.Item(74, 12).Value = 500
$ws.Cells.Item(74, 13).Value = 500
$ws.Cells.Item(74, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(74, 15).Value = "Región del Maule"
$ws.Cells.Item(74, 16).Value = 100
$ws.Cells.Item(74, 17).Value = 5
$ws.Cells.Item(74, 18).Value = "Hortaliza"
